$d = $word.ActiveDocument

# Locate the paragraph containing the old "Note:" sentence about the Notetaker
# choosing a side, by finding the distinctive old text.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Notetaker must also choose either the Employee or Employer side",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)

if (-not $found) {
    throw "Could not find target paragraph text"
}

$paraRange = $rng.Paragraphs(1).Range

$rsquo = [char]8217

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
        '<w:p>' +
            '<w:pPr><w:pStyle w:val="BlockText"/></w:pPr>' +
            '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Note:</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">If your team has</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">3 or more members</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">present, the Notetaker focuses only on facilitating and recording the debate. If only</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">2 members</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            ('<w:r><w:t xml:space="preserve">are present, the Notetaker must also choose either the Employee or Employer side. In that case, write the Notetaker' + $rsquo + 's name in both the Notetaker row AND their chosen advocate row.</w:t></w:r>') +
        '</w:p>' +
    '</w:body>' +
    '</w:document>'

$paraRange.InsertXML($newParaXml)
